# "popup delete and import data excel fix"
# Update the header/label text in the small lookup table on the first sheet:
#   A1: "name" -> "nama"
#   D2: "user" -> "users"
#   D3: "user" -> "users"
# (A2, B2, C2, A3, B3, C3, B1, C1, D1 remain unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "nama"
$ws.Range("D2").Value = "users"
$ws.Range("D3").Value = "users"
